# Final commit on 18-Jul-2020
# Rebuild the "Language"/one8-one10 sample rows into the new
# first101..104 / last201..204 / one1001..1004 / language-validation layout,
# adding a 5th data row and a new "validation" column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old sample rows so stale shared strings (Language, one8,
# one9, one10) don't linger unused in the workbook.
$ws.Range("A2:H4").ClearContents()

# ---- Header row ----
$ws.Range("A1").Value = "first name"
$ws.Range("B1").Value = "last name"
$ws.Range("C1").Value = "e-mail"
$ws.Range("D1").Value = "username"
$ws.Range("E1").Value = "password"
$ws.Range("F1").Value = "confirm password"
$ws.Range("G1").Value = "phone"
$ws.Range("H1").Value = "language"
$ws.Range("I1").Value = "validation"

# ---- Data rows ----
$ws.Range("A2").Value = "first101"
$ws.Range("B2").Value = "last201"
$ws.Range("C2").Value = "one@one.com"
$ws.Range("D2").Value = "one1001"
$ws.Range("E2").Value = "one"
$ws.Range("F2").Value = "one"
$ws.Range("G2").Value = 123
$ws.Range("H2").Value = "English"
$ws.Range("I2").Value = "An e-mail has been sent to remind you of your login and password."

$ws.Range("A3").Value = "first102"
$ws.Range("B3").Value = "last202"
$ws.Range("C3").Value = "one@one.com"
$ws.Range("D3").Value = "one1002"
$ws.Range("E3").Value = "one"
$ws.Range("F3").Value = "one"
$ws.Range("G3").Value = 123
$ws.Range("H3").Value = "English"
$ws.Range("I3").Value = "An e-mail has been sent to remind you of your login and password."

$ws.Range("A4").Value = "first103"
$ws.Range("B4").Value = "last203"
$ws.Range("C4").Value = "one@one.com"
$ws.Range("D4").Value = "one1003"
$ws.Range("E4").Value = "one"
$ws.Range("F4").Value = "one"
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = "English"
$ws.Range("I4").Value = "An e-mail has been sent to remind you of your login and password."

$ws.Range("A5").Value = "first104"
$ws.Range("B5").Value = "last204"
$ws.Range("C5").Value = "one@one.com"
$ws.Range("D5").Value = "one1004"
$ws.Range("E5").Value = "one"
$ws.Range("F5").Value = "one"
$ws.Range("G5").Value = 123
$ws.Range("H5").Value = "English"
$ws.Range("I5").Value = "An e-mail has been sent to remind you of your login and password."

# ---- Column widths: best-fit the newly populated columns ----
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(9).AutoFit()

# ---- Selection moves to E8 ----
[void]$ws.Range("E8").Select()
